$d = $word.ActiveDocument

# Update the date paragraph (first paragraph in the document)
$d.Paragraphs.Item(1).Range.Find.Execute("2022-12-26 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2022-12-27 Tuesday", 2) | Out-Null

# Update each math-problem cell in the table by position (values repeat, so
# Find/Replace across the whole doc would be ambiguous -- index by cell instead)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "74-37="
$t.Cell(1,2).Range.Text = "31-9="
$t.Cell(1,3).Range.Text = "0+20="
$t.Cell(1,4).Range.Text = "67-47="
$t.Cell(1,5).Range.Text = "77-11="
$t.Cell(2,1).Range.Text = "89-81="
$t.Cell(2,2).Range.Text = "9+57="
$t.Cell(2,3).Range.Text = "21+57="
$t.Cell(2,4).Range.Text = "22+67="
$t.Cell(2,5).Range.Text = "90+4="
$t.Cell(3,1).Range.Text = "86-10="
$t.Cell(3,2).Range.Text = "68-1="
$t.Cell(3,3).Range.Text = "15+7="
$t.Cell(3,4).Range.Text = "97-4="
$t.Cell(3,5).Range.Text = "86-0="
$t.Cell(4,1).Range.Text = "70+6="
$t.Cell(4,2).Range.Text = "44+15="
$t.Cell(4,3).Range.Text = "41-28="
$t.Cell(4,4).Range.Text = "29+51="
$t.Cell(4,5).Range.Text = "75-53="
$t.Cell(5,1).Range.Text = "0+61="
$t.Cell(5,2).Range.Text = "62-26="
$t.Cell(5,3).Range.Text = "63-38="
$t.Cell(5,4).Range.Text = "80+2="
$t.Cell(5,5).Range.Text = "9-7="
$t.Cell(6,1).Range.Text = "80+0="
$t.Cell(6,2).Range.Text = "9+60="
$t.Cell(6,3).Range.Text = "77-12="
$t.Cell(6,4).Range.Text = "86-67="
$t.Cell(6,5).Range.Text = "61-52="
$t.Cell(7,1).Range.Text = "78+9="
$t.Cell(7,2).Range.Text = "91-41="
$t.Cell(7,3).Range.Text = "55-35="
$t.Cell(7,4).Range.Text = "26+48="
$t.Cell(7,5).Range.Text = "20+8="
$t.Cell(8,1).Range.Text = "71-31="
$t.Cell(8,2).Range.Text = "8+63="
$t.Cell(8,3).Range.Text = "51+3="
$t.Cell(8,4).Range.Text = "45-20="
$t.Cell(8,5).Range.Text = "38+32="
$t.Cell(9,1).Range.Text = "7+56="
$t.Cell(9,2).Range.Text = "16+40="
$t.Cell(9,3).Range.Text = "4+74="
$t.Cell(9,4).Range.Text = "6+76="
$t.Cell(9,5).Range.Text = "38-1="
$t.Cell(10,1).Range.Text = "8+13="
$t.Cell(10,2).Range.Text = "28-4="
$t.Cell(10,3).Range.Text = "75-16="
$t.Cell(10,4).Range.Text = "34+1="
$t.Cell(10,5).Range.Text = "37-26="
$t.Cell(11,1).Range.Text = "16+32="
$t.Cell(11,2).Range.Text = "39+38="
$t.Cell(11,3).Range.Text = "41+54="
$t.Cell(11,4).Range.Text = "74-7="
$t.Cell(11,5).Range.Text = "98-76="
$t.Cell(12,1).Range.Text = "10+75="
$t.Cell(12,2).Range.Text = "59-12="
$t.Cell(12,3).Range.Text = "68-53="
$t.Cell(12,4).Range.Text = "99-86="
$t.Cell(12,5).Range.Text = "51-25="
$t.Cell(13,1).Range.Text = "20+17="
$t.Cell(13,2).Range.Text = "42+8="
$t.Cell(13,3).Range.Text = "93-77="
$t.Cell(13,4).Range.Text = "40+22="
$t.Cell(13,5).Range.Text = "97-73="
$t.Cell(14,1).Range.Text = "81-10="
$t.Cell(14,2).Range.Text = "63+9="
$t.Cell(14,3).Range.Text = "39+14="
$t.Cell(14,4).Range.Text = "0+87="
$t.Cell(14,5).Range.Text = "96-76="
$t.Cell(15,1).Range.Text = "99-89="
$t.Cell(15,2).Range.Text = "9+8="
$t.Cell(15,3).Range.Text = "85-74="
$t.Cell(15,4).Range.Text = "36-1="
$t.Cell(15,5).Range.Text = "91-88="
$t.Cell(16,1).Range.Text = "2+8="
$t.Cell(16,2).Range.Text = "30+37="
$t.Cell(16,3).Range.Text = "76-75="
$t.Cell(16,4).Range.Text = "2+54="
$t.Cell(16,5).Range.Text = "9+39="
$t.Cell(17,1).Range.Text = "94-29="
$t.Cell(17,2).Range.Text = "81+9="
$t.Cell(17,3).Range.Text = "24+62="
$t.Cell(17,4).Range.Text = "70-44="
$t.Cell(17,5).Range.Text = "6-5="
$t.Cell(18,1).Range.Text = "78-60="
$t.Cell(18,2).Range.Text = "33+24="
$t.Cell(18,3).Range.Text = "20+51="
$t.Cell(18,4).Range.Text = "85-64="
$t.Cell(18,5).Range.Text = "35+38="
$t.Cell(19,1).Range.Text = "54+13="
$t.Cell(19,2).Range.Text = "81-49="
$t.Cell(19,3).Range.Text = "73-1="
$t.Cell(19,4).Range.Text = "35+8="
$t.Cell(19,5).Range.Text = "36+26="
$t.Cell(20,1).Range.Text = "71+27="
$t.Cell(20,2).Range.Text = "71+26="
$t.Cell(20,3).Range.Text = "77+22="
$t.Cell(20,4).Range.Text = "71-4="
$t.Cell(20,5).Range.Text = "38-0="
